# Update cryptocurrency Price (D) and Volume(1h) (E) columns for rows 2-51
# to reflect the latest scrape, matching the target diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.079.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "'1.778.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.11%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'329.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.4502"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("D8").Value = "'0.3568"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'0.07455"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "'42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'1.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'21.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "'6.064"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "'1.776.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "'93.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "'0.06446"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "'0.9993"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  +2.11%  "
$ws.Range("D22").Value = "'5.803"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "'28.097.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  +1.98%  "
$ws.Range("D25").Value = "'2.123"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "'162.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'1.981.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "'2.170"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.65%  "
$ws.Range("D30").Value = "'125.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'1.108"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("D32").Value = "'5.723"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.61%  "
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("D34").Value = "'3.689"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'0.06209"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "'0.2115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").Value = "'5.009"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").Value = "'0.6333"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("D41").Value = "'1.186"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "'1.398"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").Value = "'7.932"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "'13.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("D45").Value = "'3.753"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "'0.5913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'122.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'1.964"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Value = "'1.143"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "'0.06903"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'73.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.36%  "
